$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.317.62"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.37%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.622.31"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.69%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "212.45"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.92%  "
$ws.Range("E6").Value = "  -0.02%  "
$ws.Range("E7").Value = "  +1.13%  "
$ws.Range("E8").Value = "  +1.53%  "
$ws.Range("E9").Value = "  +0.72%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "18.92"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +4.86%  "
$ws.Range("E11").Value = "  +0.85%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.848.59"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.67%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.633.17"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +2.32%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.01"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.43%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.519"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.39%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "26.325.11"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.39%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "62.51"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +4.14%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "202.97"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.58%  "
$ws.Range("E21").Value = "  +1.79%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.36"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.35%  "
$ws.Range("E23").Value = "  +0.83%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.91"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +7.24%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "143.28"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.12%  "
$ws.Range("E26").Value = "  +0.03%  "
$ws.Range("E27").Value = "  +0.42%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "15.21"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.90%  "
$ws.Range("E29").Value = "  +1.85%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0528"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +11.14%  "
$ws.Range("E31").Value = "  +0.92%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.18"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.77%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.95"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.00%  "
$ws.Range("B34").Value = "LidoDAOToken"
$ws.Range("C34").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.50"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +2.01%  "
$ws.Range("B35").Value = "HuobiToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.42"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +2.23%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.179.15"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +4.97%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0164"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.56%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.809"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +3.31%  "
$ws.Range("E39").Value = "  -0.04%  "
$ws.Range("E40").Value = "  +0.16%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.496"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.38%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.792"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.48%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.35"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +4.88%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.759.59"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.69%  "
$ws.Range("E45").Value = "  +1.13%  "
$ws.Range("E46").Value = "  +15.73%  "
$ws.Range("E47").Value = "  +1.15%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "54.01"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.39%  "
$ws.Range("E49").Value = "  +1.08%  "
$ws.Range("E50").Value = "  +0.31%  "
$ws.Range("E51").Value = "  -0.47%  "
